$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '87.225.74'
$ws.Range('E2').Value = '  -2.58%  '

$ws.Range('D3').Value = '3.143.40'
$ws.Range('E3').Value = '  -6.78%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '203.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -7.61%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '604.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.373'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -9.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.660'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.49%  '

$ws.Range('E9').Value = '  -0.09%  '

$ws.Range('D10').Value = '3.141.33'
$ws.Range('E10').Value = '  -6.80%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.531'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -12.94%  '

$ws.Range('E13').Value = '  -16.69%  '

$ws.Range('D14').Value = '3.719.26'
$ws.Range('E14').Value = '  -6.95%  '

$ws.Range('E15').Value = '  -6.71%  '

$ws.Range('D16').Value = '86.967.43'
$ws.Range('E16').Value = '  -2.76%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '31.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -13.62%  '

$ws.Range('D18').Value = '3.139.97'
$ws.Range('E18').Value = '  -6.77%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.99'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.94%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -10.51%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '411.61'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -10.30%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -12.85%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.86%  '

$ws.Range('E24').Value = '  -8.53%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.87'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -8.18%  '

$ws.Range('D26').Value = '3.308.09'
$ws.Range('E26').Value = '  -6.14%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '73.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.43%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000129'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -10.24%  '

$ws.Range('E29').Value = '  -0.16%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.162'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -20.38%  '

$ws.Range('E31').Value = '  +0.30%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '535.12'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -10.75%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.24'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -12.68%  '

$ws.Range('E34').Value = '  -18.59%  '

$ws.Range('E35').Value = '  -13.53%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.59'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -11.53%  '

$ws.Range('E37').Value = '  -8.77%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '21.78'
$ws.Range('D38').Style = 'Normal'

$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '21.77'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.51%  '

$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.16%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.97'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.01%  '

$ws.Range('E42').Value = '  +0.02%  '

$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.369'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -14.20%  '

$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -13.63%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.51'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.65%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '171.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.27%  '

$ws.Range('E47').Value = '  -7.44%  '

$ws.Range('E48').Value = '  +5.42%  '

$ws.Range('E49').Value = '  -16.00%  '

$ws.Range('E50').Value = '  -12.60%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.692'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -12.23%  '
